# Added test data reader utility
# Update the "testdata" worksheet: cell A2 header should read "runmode"
# (lowercase) instead of "Runmode".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("testdata")

$ws.Range("A2").Value = "runmode"

# Leave the selection on the edited cell, matching the saved file state.
$ws.Activate()
$ws.Range("A2").Select()
